$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A345:E345").Copy()
$ws.Range("A347:E347").PasteSpecial()
